# "Generate Report for handback" - append the handback-status row for the
# newly handed-back file adb02648-e7a4-4f4d-80f3-48f3ebfff182 to all three
# worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$fileId   = "adb02648-e7a4-4f4d-80f3-48f3ebfff182"
$handoffHash = "6a6a90ca9bbd56f40ad9a106e479d1605def8d6d"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$rowOv = 4

$wsOverview.Range("A" + $rowOv).Value = ($fileId + ".md")
$wsOverview.Range("B" + $rowOv).Value = $statusInSync
$wsOverview.Range("C" + $rowOv).Value = $statusInSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A" + $rowOv), "https://github.com/OpenLocalizationTest/oltest/blob/adb02648e7a44f4d80f348f3ebfff182/e2e/" + $fileId + ".md", "", "", $fileId + ".md") | Out-Null

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn / de-de
# Columns: Source File Name | Status | Correspond Handoff File |
#          Correspond Handoff Datetime | Target File |
#          Correspond Handback File | Correspond Handback DateTime |
#          Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$langs = @(
  @{ Name = "zh-cn"; HandoffDt = "2016-01-22 02:19:12"; HandbackDt = "2016-01-22 02:19:57" },
  @{ Name = "de-de"; HandoffDt = "2016-01-22 02:19:25"; HandbackDt = "2016-01-22 02:20:20" }
)

foreach ($lang in $langs) {
  $langName = $lang.Name
  $ws = $wb.Worksheets.Item($langName)
  $row = 4

  $handoffFileName  = $fileId + "." + $handoffHash + "." + $langName + ".xlf"
  $handbackFileName = $handoffFileName

  $ws.Range("A" + $row).Value = ($fileId + ".md")
  $ws.Range("B" + $row).Value = $statusInSync
  $ws.Range("C" + $row).Value = $handoffFileName
  $ws.Range("D" + $row).Value = $lang.HandoffDt
  $ws.Range("E" + $row).Value = ($fileId + ".md")
  $ws.Range("F" + $row).Value = $handbackFileName
  $ws.Range("G" + $row).Value = $lang.HandbackDt
  $ws.Range("H" + $row).Value = "Include"

  $ws.Range("D" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"
  $ws.Range("G" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"

  $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/adb02648e7a44f4d80f348f3ebfff182/e2e/" + $fileId + ".md"
  $targetMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest." + $langName + "/blob/" + $handoffHash + "/e2e/" + $fileId + ".md"
  $handoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $handoffHash + "/ol-handback/OpenLocalizationTestOrg/oltest." + $langName + "/xinjiang/" + $handoffFileName
  $handbackUrl = $handoffUrl

  $ws.Hyperlinks.Add($ws.Range("A" + $row), $mdUrl, "", "", $fileId + ".md") | Out-Null
  $ws.Hyperlinks.Add($ws.Range("C" + $row), $handoffUrl, "", "", $handoffFileName) | Out-Null
  $ws.Hyperlinks.Add($ws.Range("E" + $row), $targetMdUrl, "", "", $fileId + ".md") | Out-Null
  $ws.Hyperlinks.Add($ws.Range("F" + $row), $handbackUrl, "", "", $handbackFileName) | Out-Null
}
